$wb = $excel.ActiveWorkbook

# --- Sheet 1: insert the new "Dynamic_Factor" row before the existing row 24 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(24).Insert()

# Restore the merged ranges in column B (the A-column merge already grows
# correctly as part of the native row insert).
$ws1.Range("B10:B23").UnMerge()
$ws1.Range("B25:B44").UnMerge()
$ws1.Range("B10:B24").Merge()
$ws1.Range("B25:B44").Merge()

# Pick up the formatting used by the rest of the "Results" block for the
# new row's A/C cells (Insert() leaves them with a slightly different style).
$ws1.Range("A23").Copy()
$ws1.Range("A24").PasteSpecial(-4122)
$ws1.Range("C23").Copy()
$ws1.Range("C24").PasteSpecial(-4122)

$ws1.Range("C24").Value = "Dynamic_Factor"
$ws1.Range("D24").Value = 1.1
$ws1.Range("E24").Value = 1
$ws1.Range("F24").Value = 1

# --- Sheet 2: same new row, inserted before the existing row 57 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(57).Insert()

$ws2.Range("B43:B56").UnMerge()
$ws2.Range("B58:B77").UnMerge()
$ws2.Range("B43:B57").Merge()
$ws2.Range("B58:B77").Merge()

$ws2.Range("A56").Copy()
$ws2.Range("A57").PasteSpecial(-4122)
$ws2.Range("C56").Copy()
$ws2.Range("C57").PasteSpecial(-4122)

$ws2.Range("C57").Value = "Dynamic_Factor"
$ws2.Range("D57").Value = 1.1
$ws2.Range("E57").Value = 1
$ws2.Range("F57").Value = 1

Write-Output "done"
